# Commit: "Changed default delimiter character from '~' to '.'. Updated
# example files accordingly. Fixes #10"
#
# The workbook's header row used '~' as a delimiter for "array-like" field
# names (elements~C, elements~H, elements~O, vib_wavenumber~1 .. ~45).
# The new convention uses '.' instead, and the vib_wavenumber columns no
# longer carry a numeric suffix at all (they all just say
# "vib_wavenumber" - same as PyMuTT's new flattening behavior).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# elements~C / elements~H / elements~O -> elements.C / elements.H / elements.O
$ws.Range("B1").Value = "elements.C"
$ws.Range("C1").Value = "elements.H"
$ws.Range("D1").Value = "elements.O"

# vib_wavenumber~1 .. vib_wavenumber~45 (columns O..BG) -> vib_wavenumber
$vibCols = @("O","P","Q","R","S","T","U","V","W","X","Y","Z", `
             "AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM", `
             "AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW","AX","AY","AZ", `
             "BA","BB","BC","BD","BE","BF","BG")

foreach ($col in $vibCols) {
    $ws.Range($col + "1").Value = "vib_wavenumber"
}

# The author's Excel session also left the selection sitting on L1 (instead
# of the previous J16) when the file was last saved.
$ws.Range("L1").Select()
